# Updates cryptos list prices/volume deltas (and fixes the TheGraph /
# InjectiveProtocol row order) to match the latest scrape.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
# For Price cells whose new text looks like a plain number (single decimal
# point, e.g. "603.06") we force the cell to Text format first so Excel
# doesn't silently convert it to a floating point number and mangle the
# formatting (trailing zeros, exact digit count, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '70.350.96'
$ws.Cells.Item(2, 5).Value = '  +0.85%  '

$ws.Cells.Item(3, 4).Value = '3.621.90'
$ws.Cells.Item(3, 5).Value = '  +2.91%  '

$ws.Cells.Item(4, 5).Value = '  +0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '603.06'
$ws.Cells.Item(5, 5).Value = '  -0.53%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '196.21'
$ws.Cells.Item(6, 5).Value = '  -0.15%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.625'
$ws.Cells.Item(7, 5).Value = '  -0.84%  '

$ws.Cells.Item(8, 5).Value = '  +0.07%  '

$ws.Cells.Item(9, 5).Value = '  +6.77%  '

$ws.Cells.Item(10, 5).Value = '  -0.70%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '53.25'
$ws.Cells.Item(11, 5).Value = '  -0.93%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000303'
$ws.Cells.Item(12, 5).Value = '  +0.72%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '9.53'
$ws.Cells.Item(13, 5).Value = '  +0.16%  '

$ws.Cells.Item(14, 4).Value = '4.198.10'
$ws.Cells.Item(14, 5).Value = '  +2.83%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '605.62'
$ws.Cells.Item(15, 5).Value = '  +1.41%  '

$ws.Cells.Item(16, 5).Value = '  +1.08%  '

$ws.Cells.Item(17, 4).Value = '70.461.83'
$ws.Cells.Item(17, 5).Value = '  +0.77%  '

$ws.Cells.Item(18, 4).Value = '3.617.70'
$ws.Cells.Item(18, 5).Value = '  +2.41%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '19.02'
$ws.Cells.Item(19, 5).Value = '  -0.30%  '

$ws.Cells.Item(20, 5).Value = '  +1.08%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.997'
$ws.Cells.Item(21, 5).Value = '  +0.53%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '18.08'
$ws.Cells.Item(22, 5).Value = '  -1.17%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.26'
$ws.Cells.Item(23, 5).Value = '  -0.65%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '103.49'
$ws.Cells.Item(24, 5).Value = '  +1.06%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '4.60'
$ws.Cells.Item(25, 5).Value = '  -1.56%  '

$ws.Cells.Item(26, 5).Value = '  -6.61%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '10.60'
$ws.Cells.Item(27, 5).Value = '  -2.48%  '

$ws.Cells.Item(28, 5).Value = '  +0.98%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '33.77'
$ws.Cells.Item(29, 5).Value = '  +1.18%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '4.66'
$ws.Cells.Item(30, 5).Value = '  +8.10%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.20'
$ws.Cells.Item(31, 5).Value = '  +2.04%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '12.26'
$ws.Cells.Item(32, 5).Value = '  -1.43%  '

$ws.Cells.Item(33, 5).Value = '  +0.29%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '63.20'
$ws.Cells.Item(34, 5).Value = '  +0.15%  '

$ws.Cells.Item(35, 5).Value = '  +3.56%  '

$ws.Cells.Item(36, 4).Value = '3.952.72'
$ws.Cells.Item(36, 5).Value = '  +5.77%  '

$ws.Cells.Item(37, 5).Value = '  +0.16%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.07'
$ws.Cells.Item(38, 5).Value = '  -0.46%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '514.78'
$ws.Cells.Item(39, 5).Value = '  +5.38%  '

$ws.Cells.Item(40, 2).Value = 'TheGraph'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.389'
$ws.Cells.Item(40, 5).Value = '  -0.71%  '

$ws.Cells.Item(41, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '36.58'
$ws.Cells.Item(41, 5).Value = '  +0.07%  '

$ws.Cells.Item(42, 5).Value = '  -2.75%  '

$ws.Cells.Item(43, 5).Value = '  +2.61%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0460'
$ws.Cells.Item(44, 5).Value = '  +1.39%  '

$ws.Cells.Item(45, 5).Value = '  +6.80%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.88'
$ws.Cells.Item(46, 5).Value = '  +2.75%  '

$ws.Cells.Item(47, 5).Value = '  +0.30%  '

$ws.Cells.Item(48, 5).Value = '  +0.49%  '

$ws.Cells.Item(49, 5).Value = '  -0.20%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.000248'
$ws.Cells.Item(50, 5).Value = '  +0.50%  '

$ws.Cells.Item(51, 5).Value = '  +0.12%  '
